$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.603.87"
$ws.Range("E2").Value = "  +6.77%  "
$ws.Range("D3").Value = "2.394.70"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'113.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.99%  "
$ws.Range("D6").Value = "'319.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("D7").Value = "'0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "'0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").Value = "'42.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.01%  "
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").Value = "'8.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.19%  "
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "'1.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "'15.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").Value = "2.756.93"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").Value = "2.392.39"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").Value = "45.493.67"
$ws.Range("E18").Value = "  +6.09%  "
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("D21").Value = "'13.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").Value = "'74.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "'3.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.23%  "
$ws.Range("D24").Value = "'264.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("E25").Value = "  +4.55%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "'7.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.34%  "
$ws.Range("D28").Value = "'11.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.92%  "
$ws.Range("D29").Value = "'2.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("D30").Value = "'39.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.04%  "
$ws.Range("D31").Value = "'22.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "'0.0969"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.70%  "
$ws.Range("D33").Value = "'172.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("E34").Value = "  +3.56%  "
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").Value = "'4.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.06%  "
$ws.Range("E37").Value = "  +5.49%  "
$ws.Range("D38").Value = "'3.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.94%  "
$ws.Range("D39").Value = "'4.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.95%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").Value = "'1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.44%  "
$ws.Range("D42").Value = "'101.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.49%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.241"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.15%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").Value = "'13.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.07%  "
$ws.Range("D45").Value = "'71.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "'87.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.64%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'115.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.67%  "
$ws.Range("D50").Value = "'9.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.70%  "
$ws.Range("D51").Value = "1.665.76"
$ws.Range("E51").Value = "  -3.23%  "
